# Apply updated cryptocurrency price/volume data per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.482.23'
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").Value = '3.118.37'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  +0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '238.06'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.13%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '615.30'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.75%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.392'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.87%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.00%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.842'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +13.72%  '
$ws.Range("D11").Value = '3.115.10'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("E12").Value = '  -2.56%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000245'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.36%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '35.32'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").Value = '93.188.04'
$ws.Range("E15").Value = '  +1.71%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '5.44'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.19%  '
$ws.Range("D17").Value = '3.698.67'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").Value = '3.110.77'
$ws.Range("E18").Value = '  -0.87%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '3.78'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.63%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.91'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.75%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +4.26%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '443.67'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("E23").Value = '  -0.99%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.11'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -3.40%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '8.25'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +5.18%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '5.71'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.22%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '12.96'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +10.72%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '85.80'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.78%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +9.69%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.238'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("E32").Value = '  -12.82%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '9.27'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.97%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.03'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.06%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '8.04'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.13%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.161'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -9.09%  '
$ws.Range("E37").Value = '  -0.93%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.94'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -4.97%  '
$ws.Range("E39").Value = '  -1.75%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.449'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("E41").Value = '  -0.63%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '478.76'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.59%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '24.00'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +8.06%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.32'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.28%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '159.11'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("E47").Value = '  +0.53%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.87'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("E49").Value = '  -0.73%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '4.42'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.18%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '43.99'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.41%  '
